$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "dnasr281@gmail.com, System"
$newValue = "System, dnasr281@gmail.com"

$usedRows = $ws.UsedRange.Rows.Count

for ($r = 1; $r -le $usedRows; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $v = $cell.Text
    if ($v -eq $oldValue) {
        $cell.Value = $newValue
    }
}
